$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 218, shifting existing rows 218:284 down to 219:285
$ws.Rows(218).Insert()

# Populate the newly inserted row 218 with the new data record
$ws.Range("A218").Value = 4
$ws.Range("B218").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C218").Value = "Los Lagos"
$ws.Range("D218").Value = 44511
$ws.Range("E218").Value = 10
$ws.Range("F218").Value = "Fruta"
$ws.Range("G218").Value = 100102
$ws.Range("H218").Value = "Cítricos"
$ws.Range("I218").Value = 100102005
$ws.Range("J218").Value = "Naranja"
$ws.Range("K218").Value = "Lane Late"
$ws.Range("L218").Value = "Primera"
$ws.Range("M218").Value = 400
$ws.Range("N218").Value = 7500
$ws.Range("O218").Value = 8000
$ws.Range("P218").Value = 7750
$ws.Range("Q218").Value = "`$/malla 18 kilos"
$ws.Range("R218").Value = "Región de O'Higgins"
$ws.Range("S218").Value = 431
$ws.Range("T218").Value = 18
